# Add a "Prompt used" column (F) to the few-shot dataset sheet and populate it
# with "Domain_FSPrompt" for every data row, formatting the header row to
# match the new bold / center-center header style used across A1:F1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Build the new header style (bold Arial 10, centered both ways) on a
#        scratch cell, then propagate it onto A1:F1 via copy/paste-format so
#        the style table only gains the single xf we actually need (instead
#        of one per incremental property write). ---
$helper = $ws.Range("Z1")
$helper.Font.Bold = $true
$helper.HorizontalAlignment = -4108   # xlCenter
$helper.VerticalAlignment = -4108     # xlCenter
$helper.Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)  # xlPasteFormats
$helper.Clear()
$excel.CutCopyMode = $false

# --- 2. New column F: width + header + data values. ---
$ws.Columns("F").ColumnWidth = 25.7

$ws.Range("F1").Value = "Prompt used"

$promptRange = $ws.Range("F2:F41")
$promptRange.Value = "Domain_FSPrompt"

# --- 3. Update the view: scroll so row 17 is at the top and the active
#        selection sits on E36 (matches the saved workbook state). ---
$ws.Range("E36").Select()
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 1
